# "Updated agenda - Sanjay presentation."
# Slide 2 holds the IS Team Updates agenda table ("Table 1"). A new row
# for Sanjay's (TBD) 10-minute slot is inserted just above the last row
# ("Andrew" / "RA" / "5 minutes"), and the table re-flows to match.

$EMU_PER_PT = 12700.0

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the agenda table shape on the slide.
$tblShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tblShape = $candidate
    }
}

$tbl = $tblShape.Table

# Reposition the table frame first (row re-flow below fixes up heights).
$tblShape.Left = 172435 / $EMU_PER_PT
$tblShape.Top = 1445054 / $EMU_PER_PT

# Insert the new row before the final ("Andrew") row, i.e. as row 11 of 12.
$lastRowIndex = $tbl.Rows.Count
$newRow = $tbl.Rows.Add($lastRowIndex)
$newRowIndex = $lastRowIndex

$tbl.Cell($newRowIndex, 1).Shape.TextFrame.TextRange.Text = "Sanjay"
$tbl.Cell($newRowIndex, 2).Shape.TextFrame.TextRange.Text = "TBD"
$tbl.Cell($newRowIndex, 3).Shape.TextFrame.TextRange.Text = "10 minutes"

# Re-flow row heights: header/footer rows shrink slightly, and every
# interior row (including the newly-inserted one) settles at the same
# autofit height once the table grows by a row.
$rowCount = $tbl.Rows.Count
for ($i = 1; $i -le $rowCount; $i++) {
    if ($i -eq 1 -or $i -eq $rowCount) {
        $tbl.Rows.Item($i).Height = 454662 / $EMU_PER_PT
    } else {
        $tbl.Rows.Item($i).Height = 433011 / $EMU_PER_PT
    }
}

# The middle column widens by a hair as the table re-lays out.
$tbl.Columns.Item(2).Width = 7160051 / $EMU_PER_PT

# The table was re-created during this edit, bumping its shape name.
$tblShape.Name = "Table 5"
